# Natmi following Dr Hou advice
# The LR-pairs result table (Tnfsf13 -> Tnfrsf13b) was recomputed to include an
# additional "ECs" target cluster alongside the existing "M1"/"M2" targets, and
# the underlying expression/specificity statistics were refreshed accordingly.
# This expands the data block from 10 rows (5 senders x 2 targets) to 15 rows
# (5 senders x 3 targets) and updates every numeric column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs","Tnfsf13","Tnfrsf13b","ECs",3,1,1.153770666666667,3.461312,0.1294297218267158,0.1294297218267158,1,0.3333333333333333,0.4025163333333333,1.207549,0.01505239904563061,0.01505239904563061,0.4644115382542221,4.179703844287999,0.001948227821300693,0.001948227821300693),
  @("ECs","Tnfsf13","Tnfrsf13b","M1",3,1,1.153770666666667,3.461312,0.1294297218267158,0.1294297218267158,3,1,14.37655666666667,43.12967,0.5376220787283691,0.5376220787283691,16.58724936967111,149.28524432704,0.06958427609771353,0.06958427609771352),
  @("ECs","Tnfsf13","Tnfrsf13b","M2",3,1,1.153770666666667,3.461312,0.1294297218267158,0.1294297218267158,3,1,11.96193566666667,35.885807,0.4473255222260004,0.4473255222260002,13.80133048875378,124.211974398784,0.05789721790770162,0.05789721790770159),
  @("FAPs","Tnfsf13","Tnfrsf13b","ECs",3,1,0.560947,1.682841,0.06292690243138796,0.06292690243138795,1,0.3333333333333333,0.4025163333333333,1.207549,0.01505239904563061,0.01505239904563061,0.2257903296343333,2.032112966709,0.0009472008461027149,0.0009472008461027146),
  @("FAPs","Tnfsf13","Tnfrsf13b","M1",3,1,0.560947,1.682841,0.06292690243138796,0.06292690243138795,3,1,14.37655666666667,43.12967,0.5376220787283691,0.5376220787283691,8.064486332496667,72.58037699247001,0.03383089209310006,0.03383089209310005),
  @("FAPs","Tnfsf13","Tnfrsf13b","M2",3,1,0.560947,1.682841,0.06292690243138796,0.06292690243138795,3,1,11.96193566666667,35.885807,0.4473255222260004,0.4473255222260002,6.710011926409667,60.390107337687,0.02814880949218519,0.02814880949218518),
  @("M1","Tnfsf13","Tnfrsf13b","ECs",3,1,3.837574333333334,11.512723,0.4304981854736105,0.4304981854736104,1,0.3333333333333333,0.4025163333333333,1.207549,0.01505239904563061,0.01505239904563061,1.544686349547445,13.902177145927,0.006480030476168686,0.006480030476168684),
  @("M1","Tnfsf13","Tnfrsf13b","M1",3,1,3.837574333333334,11.512723,0.4304981854736105,0.4304981854736104,3,1,14.37655666666667,43.12967,0.5376220787283691,0.5376220787283691,55.17110486571223,496.5399437914101,0.2314453293631135,0.2314453293631134),
  @("M1","Tnfsf13","Tnfrsf13b","M2",3,1,3.837574333333334,11.512723,0.4304981854736105,0.4304981854736104,3,1,11.96193566666667,35.885807,0.4473255222260004,0.4473255222260002,45.90481729138456,413.1433556224611,0.1925728256343284,0.1925728256343283),
  @("M2","Tnfsf13","Tnfrsf13b","ECs",3,1,3.179188666666667,9.537566,0.3566406363494371,0.356640636349437,1,0.3333333333333333,0.4025163333333333,1.207549,0.01505239904563061,0.01505239904563061,1.279675365081556,11.517078285734,0.005368297174219362,0.00536829717421936),
  @("M2","Tnfsf13","Tnfrsf13b","M1",3,1,3.179188666666667,9.537566,0.3566406363494371,0.356640636349437,3,1,14.37655666666667,43.12967,0.5376220787283691,0.5376220787283691,45.70578602035778,411.3520741832201,0.1917378802731927,0.1917378802731927),
  @("M2","Tnfsf13","Tnfrsf13b","M2",3,1,3.179188666666667,9.537566,0.3566406363494371,0.356640636349437,3,1,11.96193566666667,35.885807,0.4473255222260004,0.4473255222260002,38.02925030286244,342.263252725762,0.159534458902025,0.159534458902025),
  @("sCs","Tnfsf13","Tnfrsf13b","ECs",3,1,0.182783,0.548349,0.02050455391884863,0.02050455391884863,1,0.3333333333333333,0.4025163333333333,1.207549,0.01505239904563061,0.01505239904563061,0.07357314295566666,0.662158286601,0.0003086427278391587,0.0003086427278391586),
  @("sCs","Tnfsf13","Tnfrsf13b","M1",3,1,0.182783,0.548349,0.02050455391884863,0.02050455391884863,3,1,14.37655666666667,43.12967,0.5376220787283691,0.5376220787283691,2.627790157203334,23.65011141483,0.01102370090124933,0.01102370090124933),
  @("sCs","Tnfsf13","Tnfrsf13b","M2",3,1,0.182783,0.548349,0.02050455391884863,0.02050455391884863,3,1,11.96193566666667,35.885807,0.4473255222260004,0.4473255222260002,2.186438486960333,19.677946382643,0.009172210289760148,0.009172210289760143)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $rowVals = $data[$i]
  for ($c = 1; $c -le $rowVals.Count; $c++) {
    $ws.Cells.Item($r, $c).Value2 = $rowVals[$c - 1]
  }
}
